# cap_data.xlsx / "Monthly 10 homes 2" sheet
# Commit: Switch from 140-149 to 70-79
# Non-solar profile changes from 140-149 to 70-79 to avoid "no calculation result" error.
#
# - Row 2 (N2:W2): sequence header values 140..149 -> 70..79
# - Rows 4-16, cols N:X: recalculated non-solar profile data + row totals
# - Active cell selection moves from Q26 to N4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monthly 10 homes 2")
$ws.Activate()

# Row 2
$ws.Cells.Item(2, 14).Value = 70  # N2: 140 -> 70
$ws.Cells.Item(2, 15).Value = 71  # O2: 141 -> 71
$ws.Cells.Item(2, 16).Value = 72  # P2: 142 -> 72
$ws.Cells.Item(2, 17).Value = 73  # Q2: 143 -> 73
$ws.Cells.Item(2, 18).Value = 74  # R2: 144 -> 74
$ws.Cells.Item(2, 19).Value = 75  # S2: 145 -> 75
$ws.Cells.Item(2, 20).Value = 76  # T2: 146 -> 76
$ws.Cells.Item(2, 21).Value = 77  # U2: 147 -> 77
$ws.Cells.Item(2, 22).Value = 78  # V2: 148 -> 78
$ws.Cells.Item(2, 23).Value = 79  # W2: 149 -> 79
# Row 4
$ws.Cells.Item(4, 14).Value = 750.30500000000006  # N4: 308.84699999999998 -> 750.30500000000006
$ws.Cells.Item(4, 15).Value = 475.37900000000002  # O4: 1200.2460000000001 -> 475.37900000000002
$ws.Cells.Item(4, 16).Value = 865.04600000000005  # P4: 427.47500000000002 -> 865.04600000000005
$ws.Cells.Item(4, 17).Value = 349.54500000000002  # Q4: 298.30799999999999 -> 349.54500000000002
$ws.Cells.Item(4, 18).Value = 612.52699999999993  # R4: 476.20899999999989 -> 612.52699999999993
$ws.Cells.Item(4, 19).Value = 763.84799999999996  # S4: 942.68299999999988 -> 763.84799999999996
$ws.Cells.Item(4, 20).Value = 632.32500000000005  # T4: 20.208000000000009 -> 632.32500000000005
$ws.Cells.Item(4, 21).Value = 544.19600000000003  # U4: 1993.1020000000005 -> 544.19600000000003
$ws.Cells.Item(4, 22).Value = 540.55799999999999  # V4: 683.02299999999991 -> 540.55799999999999
$ws.Cells.Item(4, 23).Value = 1035.6550000000002  # W4: 971.47799999999984 -> 1035.6550000000002
$ws.Cells.Item(4, 24).Value = 6569.384  # X4: 7321.5790000000006 -> 6569.384
# Row 5
$ws.Cells.Item(5, 14).Value = 675.21800000000007  # N5: 90.468000000000018 -> 675.21800000000007
$ws.Cells.Item(5, 15).Value = 403.73399999999998  # O5: 1242.4160000000002 -> 403.73399999999998
$ws.Cells.Item(5, 16).Value = 663.072  # P5: 333.59700000000004 -> 663.072
$ws.Cells.Item(5, 17).Value = 333.91  # Q5: 613.74799999999993 -> 333.91
$ws.Cells.Item(5, 18).Value = 512.00699999999995  # R5: 377.04199999999997 -> 512.00699999999995
$ws.Cells.Item(5, 19).Value = 754.32699999999977  # S5: 914.59699999999987 -> 754.32699999999977
$ws.Cells.Item(5, 20).Value = 611.33299999999997  # T5: 23.495000000000008 -> 611.33299999999997
$ws.Cells.Item(5, 21).Value = 517.75  # U5: 1581.83 -> 517.75
$ws.Cells.Item(5, 22).Value = 489.7750000000002  # V5: 872.21000000000015 -> 489.7750000000002
$ws.Cells.Item(5, 23).Value = 833.79700000000014  # W5: 920.82499999999993 -> 833.79700000000014
$ws.Cells.Item(5, 24).Value = 5794.9230000000007  # X5: 6970.2279999999992 -> 5794.9230000000007
# Row 6
$ws.Cells.Item(6, 14).Value = 577.39599999999996  # N6: 25.860000000000014 -> 577.39599999999996
$ws.Cells.Item(6, 15).Value = 340.2709999999999  # O6: 991.01699999999983 -> 340.2709999999999
$ws.Cells.Item(6, 16).Value = 589.54100000000005  # P6: 207.721 -> 589.54100000000005
$ws.Cells.Item(6, 17).Value = 283.036  # Q6: 552.56899999999985 -> 283.036
$ws.Cells.Item(6, 18).Value = 505.548  # R6: 238.55300000000003 -> 505.548
$ws.Cells.Item(6, 19).Value = 618.42800000000011  # S6: 558.07499999999993 -> 618.42800000000011
$ws.Cells.Item(6, 20).Value = 408.90699999999998  # T6: 167.68100000000001 -> 408.90699999999998
$ws.Cells.Item(6, 21).Value = 325.94499999999999  # U6: 1216.4729999999995 -> 325.94499999999999
$ws.Cells.Item(6, 22).Value = 366.62  # V6: 629.09100000000012 -> 366.62
$ws.Cells.Item(6, 23).Value = 518.25200000000007  # W6: 654.16000000000008 -> 518.25200000000007
$ws.Cells.Item(6, 24).Value = 4533.9440000000004  # X6: 5241.2 -> 4533.9440000000004
# Row 7
$ws.Cells.Item(7, 14).Value = 617.60899999999981  # N7: 100.26500000000006 -> 617.60899999999981
$ws.Cells.Item(7, 15).Value = 305.53399999999993  # O7: 961.58100000000013 -> 305.53399999999993
$ws.Cells.Item(7, 16).Value = 226.75199999999995  # P7: 332.44299999999998 -> 226.75199999999995
$ws.Cells.Item(7, 17).Value = 292.87099999999992  # Q7: 469.18200000000002 -> 292.87099999999992
$ws.Cells.Item(7, 18).Value = 492.17700000000002  # R7: 283.30100000000004 -> 492.17700000000002
$ws.Cells.Item(7, 19).Value = 893.96800000000007  # S7: 680.51999999999987 -> 893.96800000000007
$ws.Cells.Item(7, 20).Value = 328.92599999999987  # T7: 119.196 -> 328.92599999999987
$ws.Cells.Item(7, 21).Value = 297.53299999999996  # U7: 884.07700000000023 -> 297.53299999999996
$ws.Cells.Item(7, 22).Value = 373.91499999999996  # V7: 435.05400000000009 -> 373.91499999999996
$ws.Cells.Item(7, 23).Value = 362.27799999999985  # W7: 665.774 -> 362.27799999999985
$ws.Cells.Item(7, 24).Value = 4191.5630000000001  # X7: 4931.3930000000009 -> 4191.5630000000001
# Row 8
$ws.Cells.Item(8, 14).Value = 585.57899999999995  # N8: 139.77800000000002 -> 585.57899999999995
$ws.Cells.Item(8, 15).Value = 308.80800000000005  # O8: 1050.5220000000002 -> 308.80800000000005
$ws.Cells.Item(8, 16).Value = 202.89299999999997  # P8: 321.64600000000002 -> 202.89299999999997
$ws.Cells.Item(8, 17).Value = 290.65099999999995  # Q8: 593.62400000000002 -> 290.65099999999995
$ws.Cells.Item(8, 18).Value = 433.45099999999991  # R8: 323.41700000000009 -> 433.45099999999991
$ws.Cells.Item(8, 19).Value = 941.26700000000017  # S8: 751.15599999999984 -> 941.26700000000017
$ws.Cells.Item(8, 20).Value = 331.72600000000006  # T8: 130.78899999999999 -> 331.72600000000006
$ws.Cells.Item(8, 21).Value = 145.02599999999998  # U8: 673.49799999999993 -> 145.02599999999998
$ws.Cells.Item(8, 22).Value = 446.07800000000003  # V8: 524.14200000000005 -> 446.07800000000003
$ws.Cells.Item(8, 23).Value = 403.32800000000003  # W8: 639.14200000000005 -> 403.32800000000003
$ws.Cells.Item(8, 24).Value = 4088.8070000000002  # X8: 5147.7139999999999 -> 4088.8070000000002
# Row 9
$ws.Cells.Item(9, 14).Value = 706.90499999999986  # N9: 56.515999999999991 -> 706.90499999999986
$ws.Cells.Item(9, 15).Value = 326.69699999999995  # O9: 1096.4159999999997 -> 326.69699999999995
$ws.Cells.Item(9, 16).Value = 179.80199999999999  # P9: 312.05399999999997 -> 179.80199999999999
$ws.Cells.Item(9, 17).Value = 306.29199999999997  # Q9: 893.88599999999997 -> 306.29199999999997
$ws.Cells.Item(9, 18).Value = 408.25999999999993  # R9: 359.54699999999997 -> 408.25999999999993
$ws.Cells.Item(9, 19).Value = 1047.1130000000003  # S9: 887.09699999999987 -> 1047.1130000000003
$ws.Cells.Item(9, 20).Value = 337.71999999999986  # T9: 191.71199999999999 -> 337.71999999999986
$ws.Cells.Item(9, 21).Value = 154.22899999999998  # U9: 634.61099999999999 -> 154.22899999999998
$ws.Cells.Item(9, 22).Value = 513.89399999999989  # V9: 518.45800000000008 -> 513.89399999999989
$ws.Cells.Item(9, 23).Value = 472.02699999999999  # W9: 438.07399999999996 -> 472.02699999999999
$ws.Cells.Item(9, 24).Value = 4452.9389999999994  # X9: 5388.3710000000001 -> 4452.9389999999994
# Row 10
$ws.Cells.Item(10, 14).Value = 757.46800000000007  # N10: 72.271000000000001 -> 757.46800000000007
$ws.Cells.Item(10, 15).Value = 410.41600000000011  # O10: 1307.8350000000005 -> 410.41600000000011
$ws.Cells.Item(10, 16).Value = 188.42  # P10: 312.95499999999998 -> 188.42
$ws.Cells.Item(10, 17).Value = 387.02399999999994  # Q10: 970.69800000000009 -> 387.02399999999994
$ws.Cells.Item(10, 18).Value = 331.834  # R10: 404.08600000000007 -> 331.834
$ws.Cells.Item(10, 19).Value = 794.53099999999995  # S10: 948.39499999999998 -> 794.53099999999995
$ws.Cells.Item(10, 20).Value = 365.262  # T10: 239.86299999999997 -> 365.262
$ws.Cells.Item(10, 21).Value = 212.88  # U10: 711.62400000000014 -> 212.88
$ws.Cells.Item(10, 22).Value = 530.54  # V10: 557.95799999999986 -> 530.54
$ws.Cells.Item(10, 23).Value = 556.7320000000002  # W10: 509.22199999999987 -> 556.7320000000002
$ws.Cells.Item(10, 24).Value = 4535.1070000000009  # X10: 6034.9070000000002 -> 4535.1070000000009
# Row 11
$ws.Cells.Item(11, 14).Value = 561.92199999999991  # N11: 133.57500000000002 -> 561.92199999999991
$ws.Cells.Item(11, 15).Value = 377.6160000000001  # O11: 895.93299999999977 -> 377.6160000000001
$ws.Cells.Item(11, 16).Value = 180.827  # P11: 305.38900000000007 -> 180.827
$ws.Cells.Item(11, 17).Value = 280.07100000000008  # Q11: 486.78599999999994 -> 280.07100000000008
$ws.Cells.Item(11, 18).Value = 340.65400000000005  # R11: 328.66299999999995 -> 340.65400000000005
$ws.Cells.Item(11, 19).Value = 886.77499999999986  # S11: 722.3280000000002 -> 886.77499999999986
$ws.Cells.Item(11, 20).Value = 289.52900000000005  # T11: 152.81 -> 289.52900000000005
$ws.Cells.Item(11, 21).Value = 259.5150000000001  # U11: 639.41499999999996 -> 259.5150000000001
$ws.Cells.Item(11, 22).Value = 369.423  # V11: 376.78399999999999 -> 369.423
$ws.Cells.Item(11, 23).Value = 314.25399999999996  # W11: 580.29000000000008 -> 314.25399999999996
$ws.Cells.Item(11, 24).Value = 3860.5859999999993  # X11: 4621.973 -> 3860.5859999999993
# Row 12
$ws.Cells.Item(12, 14).Value = 675.10800000000006  # N12: 145.214 -> 675.10800000000006
$ws.Cells.Item(12, 15).Value = 416.53399999999993  # O12: 1028.7910000000002 -> 416.53399999999993
$ws.Cells.Item(12, 16).Value = 200.37800000000007  # P12: 337.08999999999992 -> 200.37800000000007
$ws.Cells.Item(12, 17).Value = 292.74699999999996  # Q12: 584.53300000000013 -> 292.74699999999996
$ws.Cells.Item(12, 18).Value = 354.21100000000013  # R12: 338.99300000000005 -> 354.21100000000013
$ws.Cells.Item(12, 19).Value = 986.81800000000032  # S12: 747.72100000000012 -> 986.81800000000032
$ws.Cells.Item(12, 20).Value = 324.43299999999999  # T12: 150.125 -> 324.43299999999999
$ws.Cells.Item(12, 21).Value = 293.423  # U12: 725.81299999999987 -> 293.423
$ws.Cells.Item(12, 22).Value = 492.92199999999991  # V12: 419.36500000000007 -> 492.92199999999991
$ws.Cells.Item(12, 23).Value = 393.30499999999989  # W12: 625.11400000000003 -> 393.30499999999989
$ws.Cells.Item(12, 24).Value = 4429.8790000000008  # X12: 5102.759 -> 4429.8790000000008
# Row 13
$ws.Cells.Item(13, 14).Value = 520.13099999999997  # N13: 139.90600000000001 -> 520.13099999999997
$ws.Cells.Item(13, 15).Value = 470.43400000000003  # O13: 888.79199999999992 -> 470.43400000000003
$ws.Cells.Item(13, 16).Value = 237.86199999999999  # P13: 360.97100000000006 -> 237.86199999999999
$ws.Cells.Item(13, 17).Value = 318.95699999999999  # Q13: 450.82300000000004 -> 318.95699999999999
$ws.Cells.Item(13, 18).Value = 412.42900000000003  # R13: 299.09100000000001 -> 412.42900000000003
$ws.Cells.Item(13, 19).Value = 656.51600000000019  # S13: 627.005 -> 656.51600000000019
$ws.Cells.Item(13, 20).Value = 360.36699999999996  # T13: 136.191 -> 360.36699999999996
$ws.Cells.Item(13, 21).Value = 205.851  # U13: 632.56799999999998 -> 205.851
$ws.Cells.Item(13, 22).Value = 480.70799999999997  # V13: 522.92000000000007 -> 480.70799999999997
$ws.Cells.Item(13, 23).Value = 344.57099999999991  # W13: 601.59699999999998 -> 344.57099999999991
$ws.Cells.Item(13, 24).Value = 4007.826  # X13: 4659.8639999999996 -> 4007.826
# Row 14
$ws.Cells.Item(14, 14).Value = 613.77  # N14: 171.21699999999998 -> 613.77
$ws.Cells.Item(14, 15).Value = 480.87599999999998  # O14: 989.17599999999993 -> 480.87599999999998
$ws.Cells.Item(14, 16).Value = 612.755  # P14: 373.00800000000004 -> 612.755
$ws.Cells.Item(14, 17).Value = 322.15600000000001  # Q14: 415.40100000000001 -> 322.15600000000001
$ws.Cells.Item(14, 18).Value = 575.75500000000011  # R14: 311.30799999999994 -> 575.75500000000011
$ws.Cells.Item(14, 19).Value = 635.25099999999975  # S14: 649.79500000000007 -> 635.25099999999975
$ws.Cells.Item(14, 20).Value = 516.07800000000009  # T14: 141.06000000000003 -> 516.07800000000009
$ws.Cells.Item(14, 21).Value = 282.916  # U14: 710.44699999999978 -> 282.916
$ws.Cells.Item(14, 22).Value = 572.75  # V14: 683.47599999999989 -> 572.75
$ws.Cells.Item(14, 23).Value = 554.29399999999976  # W14: 647.51700000000017 -> 554.29399999999976
$ws.Cells.Item(14, 24).Value = 5166.6009999999997  # X14: 5092.4049999999997 -> 5166.6009999999997
# Row 15
$ws.Cells.Item(15, 14).Value = 721.53199999999993  # N15: 263.52500000000003 -> 721.53199999999993
$ws.Cells.Item(15, 15).Value = 511.83499999999992  # O15: 887.87300000000005 -> 511.83499999999992
$ws.Cells.Item(15, 16).Value = 703.19900000000007  # P15: 461.42199999999985 -> 703.19900000000007
$ws.Cells.Item(15, 17).Value = 364.64300000000003  # Q15: 506.76699999999988 -> 364.64300000000003
$ws.Cells.Item(15, 18).Value = 794.42600000000004  # R15: 368.11300000000011 -> 794.42600000000004
$ws.Cells.Item(15, 19).Value = 622.15300000000002  # S15: 253.16900000000001 -> 622.15300000000002
$ws.Cells.Item(15, 20).Value = 598.69900000000007  # T15: 209.05300000000005 -> 598.69900000000007
$ws.Cells.Item(15, 21).Value = 428.50799999999992  # U15: 862.51199999999994 -> 428.50799999999992
$ws.Cells.Item(15, 22).Value = 586.85  # V15: 993.9079999999999 -> 586.85
$ws.Cells.Item(15, 23).Value = 690.64300000000026  # W15: 937.83300000000008 -> 690.64300000000026
$ws.Cells.Item(15, 24).Value = 6022.4879999999994  # X15: 5744.1750000000011 -> 6022.4879999999994
# Row 16
$ws.Cells.Item(16, 14).Value = 7762.9430000000002  # N16: 1647.442 -> 7762.9430000000002
$ws.Cells.Item(16, 15).Value = 4828.134  # O16: 12540.598 -> 4828.134
$ws.Cells.Item(16, 16).Value = 4850.5470000000005  # P16: 4085.7710000000002 -> 4850.5470000000005
$ws.Cells.Item(16, 17).Value = 3821.9029999999993  # Q16: 6836.3249999999998 -> 3821.9029999999993
$ws.Cells.Item(16, 18).Value = 5773.2790000000005  # R16: 4108.3230000000003 -> 5773.2790000000005
$ws.Cells.Item(16, 19).Value = 9600.9950000000008  # S16: 8682.5409999999993 -> 9600.9950000000008
$ws.Cells.Item(16, 20).Value = 5105.3050000000003  # T16: 1682.183 -> 5105.3050000000003
$ws.Cells.Item(16, 21).Value = 3667.7719999999995  # U16: 11265.970000000001 -> 3667.7719999999995
$ws.Cells.Item(16, 22).Value = 5764.0329999999994  # V16: 7216.3889999999992 -> 5764.0329999999994
$ws.Cells.Item(16, 23).Value = 6479.1360000000004  # W16: 8191.0259999999998 -> 6479.1360000000004
$ws.Cells.Item(16, 24).Value = 57654.046999999999  # X16: 66256.567999999999 -> 57654.046999999999

# Move the active selection to N4, matching the saved sheetView state
$ws.Range("N4").Select()

Write-Output "Applied 153 cell updates on sheet $($ws.Name)"
